# Weekly update: a new Primera/Segunda price pair is published for this
# market/product, so a new 2-row block is inserted right above the
# existing row 12 (pushing every following row down by two) and is
# populated with the same shape of data as the current newest block
# (rows 2:3), stamped with the new reporting date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 12:13 - everything from row 12 down (including
# the trailing rows) shifts down by 2, which is exactly what the diff shows.
$ws.Rows("12:13").Insert()

# Seed the new rows with the same data pattern as the most recent existing
# entry (rows 2:3), then stamp them with the new date.
$ws.Range("A2:R3").Copy($ws.Range("A12:R13"))

$ws.Range("D12").Value2 = 44882
$ws.Range("D13").Value2 = 44882
